$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.85416666666666
$ws.Range("C2").Value = 83.85416666666666
$ws.Range("D2").Value = 83.85416666666666
$ws.Range("E2").Value = 63.02083333333333
